# "start healths to 100%" - set the last station's starting health (C4 on the
# stationsLearn sheet) to 100, matching the other rows/sheets which are
# already at 100. Also reflects the author's UI navigation: they ended up
# with the "stationsLearn" sheet active/selected at C4 (and consequently
# "stationsExperiment" - previously the active tab - is no longer selected).

$wb = $excel.ActiveWorkbook

$stationsLearn = $wb.Worksheets.Item("stationsLearn")
$stationsLearn.Range("C4").Value = 100

# Make stationsLearn the active sheet and put the selection on C4, as in the
# authored workbook.
[void]$stationsLearn.Activate()
[void]$stationsLearn.Range("C4").Select()
